$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.311.64'
$ws.Range('E2').Value = '  +1.49%  '

$ws.Range('D3').Value = '3.897.05'
$ws.Range('E3').Value = '  +0.18%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.60'
$ws.Range('E5').Value = '  +9.17%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.57'
$ws.Range('E6').Value = '  -0.79%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.612'
$ws.Range('E7').Value = '  -1.35%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.719'
$ws.Range('E9').Value = '  -2.74%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.173'
$ws.Range('E10').Value = '  -1.68%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000334'
$ws.Range('E11').Value = '  -4.75%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.11'
$ws.Range('E12').Value = '  -1.86%  '

$ws.Range('D13').Value = '4.521.26'
$ws.Range('E13').Value = '  +0.36%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.25'
$ws.Range('E14').Value = '  -1.96%  '

$ws.Range('D15').Value = '3.910.93'
$ws.Range('E15').Value = '  +0.79%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.99'
$ws.Range('E16').Value = '  -1.33%  '

$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.21'
$ws.Range('E17').Value = '  +6.93%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.134'
$ws.Range('E18').Value = '  -1.39%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.78'
$ws.Range('E19').Value = '  -0.81%  '

$ws.Range('D20').Value = '69.321.17'
$ws.Range('E20').Value = '  +1.53%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '424.88'
$ws.Range('E21').Value = '  -1.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.39'
$ws.Range('E22').Value = '  -4.46%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.12'
$ws.Range('E23').Value = '  -4.18%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '87.93'
$ws.Range('E24').Value = '  -0.98%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.02'
$ws.Range('E25').Value = '  +8.82%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.45'
$ws.Range('E26').Value = '  -8.62%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.54'
$ws.Range('E27').Value = '  -3.87%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.31'
$ws.Range('E28').Value = '  -2.14%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '690.56'
$ws.Range('E29').Value = '  -4.06%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.16'
$ws.Range('E30').Value = '  -2.24%  '

$ws.Range('E31').Value = '  -2.71%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.82'
$ws.Range('E32').Value = '  -3.02%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '68.61'
$ws.Range('E33').Value = '  +10.97%  '

$ws.Range('D34').Value = '0.0₃0867'
$ws.Range('E34').Value = '  +0.11%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.431'
$ws.Range('E35').Value = '  +7.94%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.90'
$ws.Range('E36').Value = '  -2.50%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '39.84'
$ws.Range('E37').Value = '  -2.31%  '

$ws.Range('E38').Value = '  +2.44%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.04%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.00%  '

$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.27'
$ws.Range('E41').Value = '  +6.82%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0482'
$ws.Range('E42').Value = '  -2.31%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.20'
$ws.Range('E43').Value = '  +7.34%  '

$ws.Range('E44').Value = '  -6.97%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.36'
$ws.Range('E45').Value = '  +1.10%  '

$ws.Range('E46').Value = '  -1.19%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000282'
$ws.Range('E47').Value = '  +15.71%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.98'
$ws.Range('E48').Value = '  +6.69%  '

$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.751.30'
$ws.Range('E49').Value = '  +15.62%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0343'
$ws.Range('E50').Value = '  -4.74%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '145.41'
$ws.Range('E51').Value = '  +0.93%  '
